$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $val) {
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '58.165.52'
$ws.Range('E2').Value = '  -1.64%  '

Set-TextValue 'D3' '2.473.42'
$ws.Range('E3').Value = '  -1.91%  '

$ws.Range('E4').Value = '  +0.00%  '

Set-TextValue 'D5' '520.75'
$ws.Range('E5').Value = '  -2.92%  '

Set-TextValue 'D6' '132.35'
$ws.Range('E6').Value = '  -4.01%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('E8').Value = '  -1.63%  '

$ws.Range('E9').Value = '  -1.81%  '

$ws.Range('E10').Value = '  -0.83%  '

$ws.Range('E11').Value = '  +0.50%  '

$ws.Range('E12').Value = '  -1.67%  '

Set-TextValue 'D13' '2.912.10'
$ws.Range('E13').Value = '  -1.90%  '

Set-TextValue 'D14' '58.093.33'

Set-TextValue 'D15' '22.08'
$ws.Range('E15').Value = '  -4.03%  '

$ws.Range('E16').Value = '  -2.01%  '

Set-TextValue 'D17' '2.475.73'
$ws.Range('E17').Value = '  -1.69%  '

$ws.Range('E18').Value = '  -2.23%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D19' '320.87'
$ws.Range('E19').Value = '  -1.53%  '

$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D20' '4.18'
$ws.Range('E20').Value = '  -2.69%  '

$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('E22').Value = '  -3.05%  '

Set-TextValue 'D23' '64.29'
$ws.Range('E23').Value = '  -2.36%  '

$ws.Range('E24').Value = '  -3.34%  '

$ws.Range('E25').Value = '  -0.23%  '

$ws.Range('E26').Value = '  -3.53%  '

$ws.Range('E27').Value = '  -3.22%  '

Set-TextValue 'D28' '0.0₃0750'
$ws.Range('E28').Value = '  -2.60%  '

Set-TextValue 'D29' '6.37'
$ws.Range('E29').Value = '  -4.82%  '

Set-TextValue 'D30' '167.11'
$ws.Range('E30').Value = '  +2.04%  '

$ws.Range('E31').Value = '  -4.59%  '

$ws.Range('E32').Value = '  -3.55%  '

$ws.Range('E34').Value = '  -0.04%  '

Set-TextValue 'D35' '18.13'
$ws.Range('E35').Value = '  -1.89%  '

$ws.Range('E36').Value = '  -9.88%  '

Set-TextValue 'D37' '3.99'
$ws.Range('E37').Value = '  -3.20%  '

$ws.Range('E38').Value = '  -3.81%  '

Set-TextValue 'D39' '0.794'
$ws.Range('E39').Value = '  -3.01%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D40' '3.47'
$ws.Range('E40').Value = '  -4.46%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D41' '276.30'
$ws.Range('E41').Value = '  -3.45%  '

$ws.Range('E42').Value = '  -2.93%  '

$ws.Range('E43').Value = '  -1.37%  '

Set-TextValue 'D44' '126.27'
$ws.Range('E44').Value = '  -4.83%  '

$ws.Range('E45').Value = '  -2.47%  '

$ws.Range('E46').Value = '  -3.53%  '

$ws.Range('E47').Value = '  -2.87%  '

Set-TextValue 'D48' '17.16'
$ws.Range('E48').Value = '  -1.11%  '

Set-TextValue 'D49' '1.738.99'
$ws.Range('E49').Value = '  -1.53%  '

$ws.Range('E50').Value = '  -1.85%  '

$ws.Range('E51').Value = '  -1.68%  '
